# [Prototype] Add ex_PersonTimeline module
# Translate the PeopleState table headers/data from Russian/Ukrainian
# Cyrillic to their English (transliterated) equivalents.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row -----------------------------------------------------
$ws.Range("A1").Value = "FIO"
$ws.Range("B1").Value = "BirthDate"
$ws.Range("C1").Value = "City"
$ws.Range("D1").Value = "Phone"

# --- Row 2 : Ivanov ---------------------------------------------------
$ws.Range("A2").Value = "Ivanov Ivan Ivanovich"
$ws.Range("C2").Value = "Kyiv"

# --- Row 3 : Petrov -----------------------------------------------
$ws.Range("A3").Value = "Petrov Pyotr Petrovich"
$ws.Range("C3").Value = "Lviv"

# --- Row 4 : Sidorova -----------------------------------------------
$ws.Range("A4").Value = "Sidorova Anna Sergeevna"
$ws.Range("C4").Value = "Odesa"

# --- Row 5 : Kovalenko -----------------------------------------------
$ws.Range("A5").Value = "Kovalenko Maria Igorevna"
$ws.Range("C5").Value = "Kharkiv"

# --- Row 6 : Shevchenko -----------------------------------------------
$ws.Range("A6").Value = "Shevchenko Oleg Andreevich"
$ws.Range("C6").Value = "Dnipro"

# --- Row 7 : Gorbenko -----------------------------------------------
$ws.Range("A7").Value = "Gorbenko Sergey Pavlovich"
$ws.Range("C7").Value = "Kyiv"

# Rows 4-7 previously had an explicit wrapped-text row height (31.5pt)
# to fit the longer Cyrillic strings; the shorter English text fits on
# one line again, so let Excel recompute the natural row height.
$ws.Rows.Item(4).AutoFit()
$ws.Rows.Item(5).AutoFit()
$ws.Rows.Item(6).AutoFit()
$ws.Rows.Item(7).AutoFit()

# Selection moved from B4 to F4.
$ws.Range("F4").Select()
